$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '34.062.96'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -1.58%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.785.00'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -3.31%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.34%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '224.16'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -1.27%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.547'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -1.35%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.42%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '32.45'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -1.01%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.283'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -3.49%  '
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -1.73%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0930'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -0.20%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '2.045.83'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -2.88%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.818.83'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -1.11%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '10.74'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -4.11%  '
$ws.Range('B15').NumberFormat = '@'
$ws.Range('B15').Value = 'WrappedBTC'
$ws.Range('C15').NumberFormat = '@'
$ws.Range('C15').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '34.066.74'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -1.72%  '
$ws.Range('B16').NumberFormat = '@'
$ws.Range('B16').Value = 'Polygon'
$ws.Range('C16').NumberFormat = '@'
$ws.Range('C16').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.621'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -4.92%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '4.15'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -4.87%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '67.71'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -3.44%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '242.46'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -4.94%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0781'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -3.05%  '
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +0.26%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '10.67'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -5.71%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.09'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -5.34%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.10'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -3.13%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '159.48'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -1.35%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '16.26'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -4.62%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.02'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -3.15%  '
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -2.95%  '
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +0.43%  '
$ws.Range('B30').NumberFormat = '@'
$ws.Range('B30').Value = 'PancakeSwap'
$ws.Range('C30').NumberFormat = '@'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.22'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +0.71%  '
$ws.Range('B31').NumberFormat = '@'
$ws.Range('B31').Value = 'Hedera'
$ws.Range('C31').NumberFormat = '@'
$ws.Range('C31').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0514'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -3.28%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.66'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -4.60%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.50'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -3.84%  '
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -6.84%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.393.58'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -3.64%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.643'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -2.32%  '
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -2.36%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0185'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -4.39%  '
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +1.26%  '
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -0.57%  '
$ws.Range('B41').NumberFormat = '@'
$ws.Range('B41').Value = 'Aave'
$ws.Range('C41').NumberFormat = '@'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '78.42'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -5.71%  '
$ws.Range('B42').NumberFormat = '@'
$ws.Range('B42').Value = 'MXToken'
$ws.Range('C42').NumberFormat = '@'
$ws.Range('C42').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.69'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -4.55%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.909'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -7.61%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0₆0143'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +12.55%  '
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +1.71%  '
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +0.73%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '107.22'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +0.48%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '5.86'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -4.08%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.944.05'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -2.27%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '12.20'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -2.80%  '
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -0.16%  '
